$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed TPM export drops the two MuSCs-sourced edges (old rows 4 & 5)
# and only keeps the two FAPs-sourced Col9a1 -> Mag interactions.
$ws.Rows("4:5").Delete()

# Refresh the TPM-derived numeric columns (I:T) for the two remaining rows
# with the newly computed values.
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8180823333333334
$ws.Range("N2").Value = 2.454247
$ws.Range("O2").Value = 0.5115352725808422
$ws.Range("P2").Value = 0.5115352725808422
$ws.Range("Q2").Value = 0.01842266875844445
$ws.Range("R2").Value = 0.165804018826
$ws.Range("S2").Value = 0.5115352725808422
$ws.Range("T2").Value = 0.5115352725808422

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7811863333333333
$ws.Range("N3").Value = 2.343559
$ws.Range("O3").Value = 0.4884647274191578
$ws.Range("P3").Value = 0.4884647274191579
$ws.Range("Q3").Value = 0.01759179543577778
$ws.Range("R3").Value = 0.158326158922
$ws.Range("S3").Value = 0.4884647274191578
$ws.Range("T3").Value = 0.4884647274191579
